# TC27_Canine_Filter_Breed-Giant.xlsx — "corrected ICDC Breed 1-14 scripts"
#
# The FilesTab Neo4j/Cypher query (cell B4 on the "startup" sheet) is
# corrected to drop the `File Type` and `Breed` columns from its RETURN
# clause (they are not meaningful/available for the Files tab).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newQuery = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
WHERE demo.breed IN ['Giant Schnauzer'] 
OPTIONAL MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
WITH DISTINCT f, parent, c, demo, diag, s
RETURN  coalesce(f.file_name, '') AS `File Name`,
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_size, '') AS `Size`,
        coalesce(c.case_id, '') AS `Case ID`,
        coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS `Study Code`
'@

$ws.Range("B4").Value = $newQuery

# The shortened text now wraps into fewer lines, so the row shrinks.
$ws.Rows.Item(4).RowHeight = 217.5

# The selection moved from D4 to C4, scrolled so row 4 is at the top.
$ws.Range("C4").Select()
